$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "APERITIVO AL MOZART"
$ws.Range("B24").Value = "Samuele Kettmaier | A.C.DENTI"
$ws.Range("C24").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D24").Value = "Andrea Menolli | SdrumALA"
$ws.Range("E24").Value = "Giacomo  Gasparini  | Mai una gioia"
$ws.Range("F24").Value = "Alessio Zandonai | SBARX"
